$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.124.09"
$ws.Range("E2").Value = "  -3.66%  "
$ws.Range("D3").Value = "1.602.20"
$ws.Range("E3").Value = "  -2.98%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").Value = "'301.17"
$ws.Range("E6").Value = "  -2.81%  "
$ws.Range("D7").Value = "'0.3788"
$ws.Range("E7").Value = "  -2.90%  "
$ws.Range("D8").Value = "'0.3650"
$ws.Range("E8").Value = "  -4.36%  "
$ws.Range("D9").Value = "'50.06"
$ws.Range("E9").Value = "  -3.72%  "
$ws.Range("E10").Value = "  -6.82%  "
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("D12").Value = "'0.08130"
$ws.Range("E12").Value = "  -3.80%  "
$ws.Range("D13").Value = "'23.03"
$ws.Range("E13").Value = "  -3.52%  "
$ws.Range("D14").Value = "'6.600"
$ws.Range("E14").Value = "  -6.50%  "
$ws.Range("D15").Value = "'7.414"
$ws.Range("E15").Value = "  -7.30%  "
$ws.Range("E16").Value = "  -4.32%  "
$ws.Range("D17").Value = "1.603.43"
$ws.Range("E17").Value = "  -3.46%  "
$ws.Range("D18").Value = "'91.60"
$ws.Range("E18").Value = "  -2.91%  "
$ws.Range("D19").Value = "'0.06856"
$ws.Range("E19").Value = "  -2.14%  "
$ws.Range("D20").Value = "'18.28"
$ws.Range("E20").Value = "  -7.13%  "
$ws.Range("D21").Value = "'6.560"
$ws.Range("E21").Value = "  -5.95%  "
$ws.Range("B22").Value = "BitDAO"
$ws.Range("C22").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D22").Value = "'0.5568"
$ws.Range("E22").Value = "  -7.02%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "'13.01"
$ws.Range("E24").Value = "  -5.67%  "
$ws.Range("B25").Value = "WrappedBTC"
$ws.Range("C25").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D25").Value = "23.132.61"
$ws.Range("E25").Value = "  -3.60%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "'2.341"
$ws.Range("E26").Value = "  -4.10%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "'2.722"
$ws.Range("E27").Value = "  -7.28%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'21.08"
$ws.Range("E28").Value = "  -4.47%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'150.20"
$ws.Range("E29").Value = "  -1.84%  "
$ws.Range("B30").Value = "HuobiToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D30").Value = "'5.280"
$ws.Range("E30").Value = "  -2.37%  "
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").Value = "'131.91"
$ws.Range("E31").Value = "  -4.35%  "
$ws.Range("B32").Value = "WEMIXTOKEN"
$ws.Range("C32").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D32").Value = "'2.424"
$ws.Range("E32").Value = "  -4.48%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'6.845"
$ws.Range("E33").Value = "  -13.83%  "
$ws.Range("B34").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C34").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D34").Value = "1.779.99"
$ws.Range("E34").Value = "  -2.93%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.07688"
$ws.Range("E35").Value = "  -4.61%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.9460"
$ws.Range("E36").Value = "  -7.12%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02745"
$ws.Range("E37").Value = "  -6.08%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "'6.252"
$ws.Range("E38").Value = "  -7.18%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "'0.2544"
$ws.Range("E39").Value = "  -4.89%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").Value = "'0.08932"
$ws.Range("E40").Value = "  -1.80%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'10.02"
$ws.Range("E41").Value = "  -6.65%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.387"
$ws.Range("E42").Value = "  -2.36%  "
$ws.Range("D43").Value = "'0.7109"
$ws.Range("E43").Value = "  -6.53%  "
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").Value = "'12.71"
$ws.Range("E44").Value = "  -5.16%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'15.45"
$ws.Range("E45").Value = "  -4.89%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.6621"
$ws.Range("E46").Value = "  -4.96%  "
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "'1.001"
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'2.300"
$ws.Range("E48").Value = "  -6.39%  "
$ws.Range("B49").Value = "PancakeSwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D49").Value = "'3.973"
$ws.Range("E49").Value = "  -3.05%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'132.09"
$ws.Range("E50").Value = "  -2.11%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.07937"
$ws.Range("E51").Value = "  -4.61%  "
